# Apply spell-checker re-run results: rewrite B/C columns for rows 2-26
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 2).Value = "previnem, reactions, recommended, saudáveis​​, anvisa, triplice, obukhanych, vaxtruth, vaccinated, transgênicos, unvaccinated, mandates, empurradores, ​​e, nauseam, vaers, reorientar, studies, ogms, sites, –, tetyana, healthier, defending, ’, vactruthcom, gardasil, diseases, sgb, pharma, parents, says, reasons, report, gms, comprovadamente, ingredients, updated, “, eficazmente, sanevax, vaccines, higher, —, imunocomprometidas, ​​em, origins, rentáveis, injury, vactruth, related, firms"
$ws.Cells.Item(2, 3).Value = 0.03649122807017544

# Row 3
$ws.Cells.Item(3, 2).Value = "diretaço, bronchiale, adaptabilidade, phds, sites, coróide, subsequentemente, colecisti, linfomas, jfj, candidíases, epatocarcinoma, tumore, ingerível, deleto, adequadamente, turmores, trattati, pediátrica, nahco3, alla, è, macrobióticos, carcinosi, peritoneale, estomago, alcalinizando, vescica, charcot, polmonari, tullio, “, metodicamente, piccoli, prostata, poderosíssima, pazes, metastasipolmonari, pessoalmente, 360°, reversões, restringe, cânceres, polmoni, dias3, pecezinho, •, diffuso, simoncini, gotejador, links, 90°, intertítulo, cerebrali, absurdamente, abrasões, casi, lincados, subministrar, tanti, cérvico, entopem, 5todos, semanas4, espetaculares, alcuni, coliciste, dellintestino, vezes1, mês1, oncologista, ½, inalador, legendado, tumori, radiografias, palpebra, midollare, oncologistas, 36a"
$ws.Cells.Item(3, 3).Value = 0.01900688999762414

# Row 4
$ws.Cells.Item(4, 2).Value = "refrescante, 16º, cm2, “, simplemente"
$ws.Cells.Item(4, 3).Value = 0.01436781609195402

# Row 5
$ws.Cells.Item(5, 2).Value = "adoçante, blaylock, 30º, arckle, entitulado, atenciosamente, câimbras, mancy, zerocal, dopamina, riccio, hj, raffaele, lobbies"
$ws.Cells.Item(5, 3).Value = 0.01690821256038647

# Row 6
$ws.Cells.Item(6, 2).Value = "possivel"
$ws.Cells.Item(6, 3).Value = 0.0101010101010101

# Row 7
$ws.Cells.Item(7, 2).Value = "cairam, lugarquando"
$ws.Cells.Item(7, 3).Value = 0.006557377049180328

# Row 8
$ws.Cells.Item(8, 2).Value = "reage, vómitos"
$ws.Cells.Item(8, 3).Value = 0.003273322422258593

# Row 9
$ws.Cells.Item(9, 2).Value = "infectologista, –, nossso, tamiflu, famíliares, h1n1"
$ws.Cells.Item(9, 3).Value = 0.03821656050955414

# Row 10
$ws.Cells.Item(10, 2).Value = "fenofinol, socbrasde, ambev, “, voliteral, skol, fleury, almeido"
$ws.Cells.Item(10, 3).Value = 0.03333333333333333

# Row 11
$ws.Cells.Item(11, 2).Value = "fenofinol, fanta, socbrasde, “, voliteral, proprios, fleury, almeido"
$ws.Cells.Item(11, 3).Value = 0.02797202797202797

# Row 12
$ws.Cells.Item(12, 2).Value = "osh, polifenol, guruprasad, handphone, reddy"
$ws.Cells.Item(12, 3).Value = 0.008665511265164644

# Row 13
$ws.Cells.Item(13, 2).Value = "laetril, mirtilos, “, alforjón, b17, amígdalina, 5sementes"
$ws.Cells.Item(13, 3).Value = 0.01107594936708861

# Row 14
$ws.Cells.Item(14, 2).Value = "11h, 16h"
$ws.Cells.Item(14, 3).Value = 0.01

# Row 15
$ws.Cells.Item(15, 2).Value = "lpki"
$ws.Cells.Item(15, 3).Value = 0.003759398496240601

# Row 16
$ws.Cells.Item(16, 2).Value = $null
$ws.Cells.Item(16, 3).Value = 0.0

# Row 17
$ws.Cells.Item(17, 2).Value = "sacarovictus, prontamente, superdivertida, –, ‘, aconchegante, ’, cevabacillus, ativus, humildemente, libre, contactei"
$ws.Cells.Item(17, 3).Value = 0.0273972602739726

# Row 18
$ws.Cells.Item(18, 2).Value = "50cc, desconfortáveis, 80cc"
$ws.Cells.Item(18, 3).Value = 0.01477832512315271

# Row 19
$ws.Cells.Item(19, 2).Value = "informaçao, podera, transmissao, afectados, japao, pergosa, adiquire, nao, estao, medicos"
$ws.Cells.Item(19, 3).Value = 0.136986301369863

# Row 20
$ws.Cells.Item(20, 2).Value = $null
$ws.Cells.Item(20, 3).Value = 0.0

# Row 21
$ws.Cells.Item(21, 2).Value = "começõu, globalista, wurhan, rothschilds, rockfellers, lives, globalistas, alcool, jinping, totall"
$ws.Cells.Item(21, 3).Value = 0.05747126436781609

# Row 22
$ws.Cells.Item(22, 2).Value = "19h25, ​​e"
$ws.Cells.Item(22, 3).Value = 0.004987531172069825

# Row 23
$ws.Cells.Item(23, 2).Value = "reencaminhado"
$ws.Cells.Item(23, 3).Value = 0.003021148036253776

# Row 24
$ws.Cells.Item(24, 2).Value = "ademola, familiares"
$ws.Cells.Item(24, 3).Value = 0.01162790697674419

# Row 25
$ws.Cells.Item(25, 2).Value = "inalador"
$ws.Cells.Item(25, 3).Value = 0.004504504504504504

# Row 26
$ws.Cells.Item(26, 2).Value = "huiren, cistos, polifenóis, “, familiares, anticâncer"
$ws.Cells.Item(26, 3).Value = 0.02608695652173913
